# FINFLUX-2698 Correcting Overdue failed scenarios
#
# On the "Transactions" sheet:
#   - Column A (loan/transaction #) values are renumbered (each decreased
#     by 345).
#   - The amount that used to sit in the "Penalties" column (I) is moved
#     into the "Amount" column (E) for every data row (2-22), leaving the
#     Penalties column blank. For the big round-number "Disbursement" /
#     "Repayment" rows (7, 10, 15, 22) the special #,##0 number format
#     that used to live on column I travels with the value to columns
#     E:H as well.
#   - The active selection on the sheet moves to E20.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")
$ws.Activate()

# New column-A numbers, row by row (row => new value)
$newIds = @{
    2  = 555
    3  = 554
    4  = 553
    5  = 552
    6  = 551
    7  = 546
    8  = 540
    9  = 539
    10 = 538
    11 = 545
    12 = 524
    13 = 523
    14 = 522
    15 = 513
    16 = 504
    17 = 495
    18 = 494
    19 = 493
    20 = 482
    21 = 481
    22 = 480
}

# Rows whose Penalties (I) cell carries the special "#,##0" number format
# (style shared across E:I) that must travel with the value.
$specialFormatRows = @(7, 10, 15, 22)

foreach ($row in 2..22) {

    # --- Column A renumbering ---
    $ws.Cells.Item($row, 1).Value2 = $newIds[$row]

    # --- Move the Penalties value (I) into Amount (E) ---
    $penaltyCell = $ws.Cells.Item($row, 9)   # column I
    $amountCell  = $ws.Cells.Item($row, 5)   # column E
    $amountValue = $penaltyCell.Value2

    if ($specialFormatRows -contains $row) {
        # Carry the distinctive "#,##0" number format that used to live
        # only on I across the whole E:I block (I itself already has it,
        # and keeps it - only its value is cleared).
        $ws.Range($ws.Cells.Item($row, 9), $ws.Cells.Item($row, 9)).Copy() | Out-Null
        $destRange = $ws.Range($ws.Cells.Item($row, 5), $ws.Cells.Item($row, 8))
        $destRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

        $amountCell.Value2 = $amountValue
        $penaltyCell.Value2 = $null
    }
    else {
        $amountCell.Value2 = $amountValue
        $penaltyCell.Value2 = $null
    }
}

$excel.CutCopyMode = $false

# --- Sheet selection / view state ---
$ws.Range("E20").Select()

$win = $excel.ActiveWindow
try { $win.ScrollRow = 19 } catch {}
try { $win.ScrollColumn = 1 } catch {}
try { $win.FirstSheet = 4 } catch {}
